# The "Id" column (integer primary key, column A) was causing int-parsing
# problems downstream, so remove it entirely. This shifts UniqueId/Name/Email
# left into columns A/B/C and shrinks the used range from A1:D6 to A1:C6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").EntireColumn.Delete()

# Move/restore the active selection to sit just below the now-shorter table.
$ws.Range("A6").Select()
